# Update cryptos list values (Price and Volume(1h) columns) per scheduled GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.579.51'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.116.02'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("D4").Value = '''1.011'
$ws.Range("E4").Value = '  +0.85%  '
$ws.Range("D5").Value = '''336.78'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("D7").Value = '''0.5245'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("D8").Value = '''0.4562'
$ws.Range("E8").Value = '  +4.12%  '
$ws.Range("D9").Value = '''54.29'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").Value = '''0.09111'
$ws.Range("E10").Value = '  +2.08%  '
$ws.Range("D11").Value = '''1.174'
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = '''24.62'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '2.122.20'
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").Value = '''6.870'
$ws.Range("E14").Value = '  +2.60%  '
$ws.Range("D15").Value = '''8.117'
$ws.Range("E15").Value = '  +5.59%  '
$ws.Range("D16").Value = '''0.00001174'
$ws.Range("E16").Value = '  +4.56%  '
$ws.Range("D17").Value = '''97.09'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '''1.011'
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("D20").Value = '''19.44'
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("D22").Value = '''6.312'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = '30.647.71'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").Value = '''12.86'
$ws.Range("E24").Value = '  +4.74%  '
$ws.Range("D25").Value = '''2.360'
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("D26").Value = '2.365.34'
$ws.Range("E26").Value = '  +1.53%  '
$ws.Range("D27").Value = '''22.36'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").Value = '''163.91'
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").Value = '''2.550'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").Value = '''134.18'
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").Value = '''1.213'
$ws.Range("E31").Value = '  +2.48%  '
$ws.Range("D32").Value = '''0.1073'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").Value = '''1.658'
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").Value = '''6.376'
$ws.Range("E34").Value = '  +3.43%  '
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").Value = '''10.58'
$ws.Range("E36").Value = '  +5.60%  '
$ws.Range("D37").Value = '''5.912'
$ws.Range("E37").Value = '  +7.87%  '
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("D39").Value = '''0.06838'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").Value = '''0.2332'
$ws.Range("E40").Value = '  +3.34%  '
$ws.Range("D41").Value = '''12.63'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").Value = '''0.6889'
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '''1.259'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").Value = '''14.95'
$ws.Range("E44").Value = '  +6.75%  '
$ws.Range("E45").Value = '  +1.91%  '
$ws.Range("D46").Value = '''2.321'
$ws.Range("E46").Value = '  +5.66%  '
$ws.Range("D47").Value = '''0.00000000368'
$ws.Range("E47").Value = '  +22.54%  '
$ws.Range("D48").Value = '''3.689'
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("D49").Value = '''1.255'
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("D50").Value = '''83.56'
$ws.Range("E50").Value = '  +2.15%  '
$ws.Range("D51").Value = '''0.3350'
$ws.Range("E51").Value = '  +12.19%  '
